$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Memory Usage (bytes)"

$ws.Range("C2").Value = 16.08705520629883
$ws.Range("C3").Value = 15.69008827209473
$ws.Range("C4").Value = 16.47377014160156
$ws.Range("C5").Value = 17.90094375610352
$ws.Range("C6").Value = 17.98701286315918
